# Thesis "progress" tracker: append the day's writing-progress row (4/16/2020 data
# point for 12/... actually 2020-11-10, serial 44136) to the "writing" table, and
# keep the dashboard chart's series ranges pointing at the now-larger range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("writing")

# --- 1. Grow Table1 by one row (keeps ref/autoFilter/dimension in sync) ---
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

# --- 2. Copy the date formatting down from A7 into the new A8 cell, then set values ---
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A8").Value = 44136
$ws.Range("B8").Value = 239
$ws.Range("C8").Value = 87
$ws.Range("D8").Value = 528
$ws.Range("E8").Value = 6658
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 6
$ws.Range("H8").Value = 5
$ws.Range("I8").Value = 5
$ws.Range("J8").Formula = "=SUM(B8:I8)"
$ws.Range("K8").Value = 621

# --- 3. Point the dashboard chart's two series at the new $A$2:$A$8 / $J$ / $K$ ranges ---
$dash = $wb.Worksheets.Item("dashboard")
$chart = $dash.ChartObjects().Item(1).Chart

$daily = $chart.SeriesCollection().Item(1)
$daily.Formula = "=SERIES(writing!`$K`$1,writing!`$A`$2:`$A`$8,writing!`$K`$2:`$K`$8,1)"

$total = $chart.SeriesCollection().Item(2)
$total.Formula = "=SERIES(writing!`$J`$1,writing!`$A`$2:`$A`$8,writing!`$J`$2:`$J`$8,2)"

# --- 4. Leave the cursor where the author left it after typing the new row ---
$ws.Activate()
$ws.Range("B9").Select()
